# Swap the contents of two full rows (only for the given column letters),
# leaving row numbers/positions untouched but exchanging all data between
# them. Uses Range.Copy so the original cell type (number/text/bool/date)
# is preserved exactly. Copying an empty source cell onto a range is a
# no-op (it does not clear the destination), so empty/blank source cells
# are instead handled with ClearContents.

function Copy-CellValue {
    param($ws, $srcAddr, $dstAddr)

    $srcRange = $ws.Range($srcAddr)
    $v = $srcRange.Value2
    $isEmpty = ($null -eq $v) -or (($v -is [string]) -and ($v.Length -eq 0))

    if ($isEmpty) {
        $ws.Range($dstAddr).ClearContents()
    } else {
        $srcRange.Copy($ws.Range($dstAddr))
    }
}

function Swap-RowData {
    param($ws, $row1, $row2, $cols, $scratchAddr)

    foreach ($c in $cols) {
        $addr1 = $c + $row1
        $addr2 = $c + $row2

        Copy-CellValue $ws $addr1 $scratchAddr
        Copy-CellValue $ws $addr2 $addr1
        Copy-CellValue $ws $scratchAddr $addr2
    }

    $ws.Range($scratchAddr).ClearContents()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All columns that are populated anywhere in the affected rows.
$cols = @("A","B","C","D","E","F","G","H","I","J","K","P","Q","R","S","T","U","V","W","Y","Z","AA","AB","AD","AE","AG","AT","AW","AX","AY")

# Far-away scratch cell used as temporary swap storage; cleared afterwards.
$scratchAddr = "ZZ1"

Swap-RowData $ws 71 72 $cols $scratchAddr
Swap-RowData $ws 97 98 $cols $scratchAddr
Swap-RowData $ws 100 101 $cols $scratchAddr
